$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $value) {
    $range = $sheet.Range($addr)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

Set-CellText $ws 'D2' '29.783.33'
Set-CellText $ws 'E2' '  +3.20%  '
Set-CellText $ws 'D3' '1.852.57'
Set-CellText $ws 'E3' '  +2.43%  '
Set-CellText $ws 'D4' '0.9989'
Set-CellText $ws 'E4' '  +0.03%  '
Set-CellText $ws 'D5' '245.08'
Set-CellText $ws 'E5' '  +2.22%  '
Set-CellText $ws 'D6' '0.6399'
Set-CellText $ws 'E6' '  +5.73%  '
Set-CellText $ws 'D7' '0.9998'
Set-CellText $ws 'E7' '  +0.24%  '
Set-CellText $ws 'B8' 'Cardano'
Set-CellText $ws 'C8' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-CellText $ws 'D8' '0.3011'
Set-CellText $ws 'E8' '  +5.43%  '
Set-CellText $ws 'B9' 'Dogecoin'
Set-CellText $ws 'C9' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-CellText $ws 'D9' '0.07517'
Set-CellText $ws 'E9' '  +3.63%  '
Set-CellText $ws 'D10' '24.15'
Set-CellText $ws 'E10' '  +6.42%  '
Set-CellText $ws 'D11' '0.07682'
Set-CellText $ws 'E11' '  +0.93%  '
Set-CellText $ws 'D12' '1.856.39'
Set-CellText $ws 'E12' '  +3.11%  '
Set-CellText $ws 'D13' '5.069'
Set-CellText $ws 'E13' '  +3.37%  '
Set-CellText $ws 'D14' '0.6899'
Set-CellText $ws 'E14' '  +5.49%  '
Set-CellText $ws 'D15' '84.60'
Set-CellText $ws 'E15' '  +4.89%  '
Set-CellText $ws 'D16' '0.000009622'
Set-CellText $ws 'E16' '  +8.21%  '
Set-CellText $ws 'D17' '6.098'
Set-CellText $ws 'E17' '  +4.85%  '
Set-CellText $ws 'D18' '29.750.60'
Set-CellText $ws 'E18' '  +3.24%  '
Set-CellText $ws 'D19' '2.092.30'
Set-CellText $ws 'E19' '  +2.72%  '
Set-CellText $ws 'D20' '240.00'
Set-CellText $ws 'E20' '  +1.78%  '
Set-CellText $ws 'D21' '12.65'
Set-CellText $ws 'E21' '  +2.64%  '
Set-CellText $ws 'D22' '0.9998'
Set-CellText $ws 'E22' '  +0.20%  '
Set-CellText $ws 'D23' '7.364'
Set-CellText $ws 'E23' '  +4.29%  '
Set-CellText $ws 'D24' '0.9992'
Set-CellText $ws 'E24' '  +0.05%  '
Set-CellText $ws 'D25' '159.74'
Set-CellText $ws 'E25' '  +1.51%  '
Set-CellText $ws 'D26' '0.1423'
Set-CellText $ws 'E26' '  +2.08%  '
Set-CellText $ws 'D27' '8.554'
Set-CellText $ws 'E27' '  +2.50%  '
Set-CellText $ws 'D28' '17.96'
Set-CellText $ws 'E28' '  +2.80%  '
Set-CellText $ws 'D29' '1.502'
Set-CellText $ws 'E29' '  +2.16%  '
Set-CellText $ws 'D30' '0.06068'
Set-CellText $ws 'E30' '  +9.32%  '
Set-CellText $ws 'D31' '1.263'
Set-CellText $ws 'E31' '  +5.22%  '
Set-CellText $ws 'D32' '4.149'
Set-CellText $ws 'E32' '  +2.40%  '
Set-CellText $ws 'D33' '4.142'
Set-CellText $ws 'E33' '  +2.11%  '
Set-CellText $ws 'D34' '1.871'
Set-CellText $ws 'E34' '  +3.72%  '
Set-CellText $ws 'B35' 'ARBITRUM'
Set-CellText $ws 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-CellText $ws 'D35' '1.153'
Set-CellText $ws 'E35' '  +2.79%  '
Set-CellText $ws 'B36' 'ImmutableX'
Set-CellText $ws 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws 'D36' '0.7347'
Set-CellText $ws 'E36' '  +0.90%  '
Set-CellText $ws 'D37' '2.611'
Set-CellText $ws 'E37' '  -0.13%  '
Set-CellText $ws 'D38' '2.859'
Set-CellText $ws 'E38' '  +2.41%  '
Set-CellText $ws 'D39' '1.228.55'
Set-CellText $ws 'E39' '  +3.61%  '
Set-CellText $ws 'D40' '0.01784'
Set-CellText $ws 'E40' '  +2.56%  '
Set-CellText $ws 'D41' '6.374'
Set-CellText $ws 'E41' '  +1.00%  '
Set-CellText $ws 'D42' '0.9226'
Set-CellText $ws 'E42' '  +4.64%  '
Set-CellText $ws 'E43' '  +0.54%  '
Set-CellText $ws 'D44' '2.014.19'
Set-CellText $ws 'E44' '  +3.71%  '
Set-CellText $ws 'D45' '102.38'
Set-CellText $ws 'D46' '66.61'
Set-CellText $ws 'E46' '  +4.39%  '
Set-CellText $ws 'D47' '0.00000000121'
Set-CellText $ws 'E47' '  -0.22%  '
Set-CellText $ws 'D48' '0.5084'
Set-CellText $ws 'E48' '  +0.46%  '
Set-CellText $ws 'D49' '9.299'
Set-CellText $ws 'E49' '  +3.86%  '
Set-CellText $ws 'D50' '0.4090'
Set-CellText $ws 'E50' '  +3.60%  '
Set-CellText $ws 'D51' '0.1146'
Set-CellText $ws 'E51' '  +4.88%  '
